# Generate Report for Handback
#
# Refresh the handback-status report after a new handback round-trip for
# the "281112d5-af5a-4b77-ac95-317023dcf2e9.md" source file: its
# zh-cn and de-de handoff/handback timestamps move forward, and the
# Overview sheet's "Latest HO Xliff Generate Date" for that file is
# recomputed as the max of those timestamps. The other file
# ("72dfbeb8-...md") had no new activity, so its rows stay as-is.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# zh-cn: Correspond Handoff Datetime (H) / Correspond Handback DateTime (K)
# for the 281112d5-... row (row 2).
$zhcn.Range("H2").Value = "2016-08-15 22:45:18"
$zhcn.Range("K2").Value = "2016-08-15 22:45:44"

# de-de: same two columns, same row.
$dede.Range("H2").Value = "2016-08-15 22:45:24"
$dede.Range("K2").Value = "2016-08-15 22:45:53"

# Overview: Latest HO Xliff Generate Date (G) for the 281112d5-... row
# (row 2) becomes the newest of the timestamps just written above.
$overview.Range("G2").Value = "2016-08-15 22:45:24"
